$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-02 Friday" "2024-08-03 Saturday"

Replace-Text "510÷7=72, 6" "104÷4=26, 0"
Replace-Text "134÷8=16, 6" "518÷4=129, 2"
Replace-Text "930÷5=186, 0" "351÷3=117, 0"
Replace-Text "205÷3=68, 1" "151÷6=25, 1"
Replace-Text "109÷4=27, 1" "843÷2=421, 1"

Replace-Text "753÷6=125, 3" "614÷3=204, 2"
Replace-Text "848÷3=282, 2" "351÷3=117, 0"
Replace-Text "340÷7=48, 4" "719÷8=89, 7"
Replace-Text "369÷6=61, 3" "746÷4=186, 2"
Replace-Text "950÷8=118, 6" "610÷4=152, 2"

Replace-Text "793÷7=113, 2" "325÷7=46, 3"
Replace-Text "538÷3=179, 1" "827÷2=413, 1"
Replace-Text "928÷6=154, 4" "383÷4=95, 3"
Replace-Text "663÷6=110, 3" "559÷3=186, 1"
Replace-Text "415÷2=207, 1" "226÷2=113, 0"

Replace-Text "712÷4=178, 0" "565÷5=113, 0"
Replace-Text "744÷8=93, 0" "265÷9=29, 4"
Replace-Text "536÷9=59, 5" "834÷2=417, 0"
Replace-Text "595÷2=297, 1" "365÷5=73, 0"
Replace-Text "141÷6=23, 3" "642÷4=160, 2"

Replace-Text "711÷7=101, 4" "810÷3=270, 0"
Replace-Text "879÷7=125, 4" "890÷5=178, 0"
Replace-Text "855÷2=427, 1" "750÷4=187, 2"
Replace-Text "245÷7=35, 0" "584÷4=146, 0"
Replace-Text "742÷4=185, 2" "987÷8=123, 3"
